# ============================================================================
# Applies the GRAFS-e example workbook update:
#   - "crops" sheet: insert 3 new columns (Spreading Rate / Seed input / Area)
#     before the existing "Carbon Mechanisation..." block, and append a new
#     "Forage crop" row.
#   - "excretion" sheet: append "Methanization power" + "Nitrogen Content"
#     columns, and bump two grasslands-excretion values.
#   - "prod" sheet: append "Methanization power" + "Production (kton)"
#     columns, and append a new "Forage crop" row.
#   - "global" sheet: append 6 new methanizer/green-waste parameter rows and
#     make it the active sheet.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "crops" (sheet1)
# ----------------------------------------------------------------------------
$crops = $wb.Worksheets.Item("crops")

# Insert 3 fresh columns at J:L - this shifts the old J..M ("Carbon
# Mechanisation Intensity", "Root Humification Coefficient", "Residue
# Humification Coefficient", "Surface Root Production") block to M..P,
# carrying over their values/styles untouched.
$crops.Range("J1:L1").EntireColumn.Insert()

# Approximate column widths (chars) for the new/shifted columns.
$crops.Columns.Item(10).ColumnWidth = 15.75
$crops.Range("K1:L1").EntireColumn.ColumnWidth = 24.584

# New headers for the inserted columns.
$crops.Range("J1").Value = "Spreading Rate (%)"
$crops.Range("K1").Value = "Seed input (kt seeds/kt Ymax)"
$crops.Range("L1").Value = "Area (ha)"

# New columns are all zero for every existing crop row (2-8); copy the
# number style (s="4") from a neighbouring numeric cell so no spurious new
# style entries get created.
$crops.Range("J2:L8").Value = 0
$crops.Range("E2").Copy()
$crops.Range("J2:L8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "Forage crop" row (row 9).
$crops.Range("A9").Value = "Forage crop"
$crops.Range("B9").Value = "Forage"
$crops.Range("C9").Value = "natural meadows"
$crops.Range("D9").Value = 0
$crops.Range("E9:P9").Value = 0

$crops.Range("B2").Copy()
$crops.Range("B9").PasteSpecial(-4122)
$crops.Range("C2").Copy()
$crops.Range("C9").PasteSpecial(-4122)
$crops.Range("D6").Copy()
$crops.Range("D9").PasteSpecial(-4122)
$crops.Range("E2").Copy()
$crops.Range("E9:P9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# View state: selection moves to N2, top-left scrolls to column E, and this
# sheet is no longer the active tab (the workbook now opens on "global").
$crops.Range("N2").Select()

# ----------------------------------------------------------------------------
# Sheet "excretion" (sheet4)
# ----------------------------------------------------------------------------
$excretion = $wb.Worksheets.Item("excretion")

$excretion.Columns.Item(4).ColumnWidth = 15.251
$excretion.Columns.Item(10).ColumnWidth = 17.75
$excretion.Columns.Item(11).ColumnWidth = 18.917

$excretion.Range("J1").Value = "Methanization power (MWh/tFW)"
$excretion.Range("K1").Value = "Nitrogen Content (%)"

$excretion.Range("J2").Value = 0.23
$excretion.Range("K2").Value = 0.5

$excretion.Range("J3").Value = 0.28999999999999998
$excretion.Range("K3").Value = 0.3

$excretion.Range("J4").Value = 0
$excretion.Range("K4").Value = 0

$excretion.Range("I5").Value = 20
$excretion.Range("J5").Value = 0.54
$excretion.Range("K5").Value = 0.8

$excretion.Range("I6").Value = 10
$excretion.Range("J6").Value = 0.21
$excretion.Range("K6").Value = 0.5

$excretion.Range("J7").Value = 0
$excretion.Range("K7").Value = 0

# ----------------------------------------------------------------------------
# Sheet "prod" (sheet5)
# ----------------------------------------------------------------------------
$prod = $wb.Worksheets.Item("prod")

$prod.Range("I1").Value = "Methanization power (MWh/tFW)"
$prod.Range("J1").Value = "Production (kton)"

$prod.Range("I2").Value = 0.78
$prod.Range("I3").Value = 0.78
$prod.Range("I4").Value = 0.78
$prod.Range("I5").Value = 0.78
$prod.Range("I6").Value = 0.78
$prod.Range("I7").Value = 0.25
$prod.Range("I8").Value = 0.25
$prod.Range("I9").Value = 0
$prod.Range("I10").Value = 0
$prod.Range("I11").Value = 0
$prod.Range("I12").Value = 0

# New "Forage crop" row (row 13).
$prod.Range("A13").Value = "Forage"
$prod.Range("B13").Value = "Forage crop"
$prod.Range("C13").Value = "plant"
$prod.Range("D13").Value = "forage"
$prod.Range("E13").Value = 2
$prod.Range("F13").Value = 45
$prod.Range("G13").Value = 0
$prod.Range("H13").Value = 0
$prod.Range("I13").Value = 0.25
$prod.Range("J13").Value = 0

$prod.Range("H14").Select()

# ----------------------------------------------------------------------------
# Sheet "global" (sheet6)
# ----------------------------------------------------------------------------
$global = $wb.Worksheets.Item("global")

$global.Range("A7").Value = "Methanizer Energy Production (GWh)"
$global.Range("B7").Value = 1000

$global.Range("A8").Value = "Weight methanizer production"
$global.Range("B8").Value = 1

$global.Range("A9").Value = "Weight methanizer inputs"
$global.Range("B9").Value = 1

$global.Range("A10").Value = "Green waste methanization power (MWh/ktN)"
$global.Range("B10").Value = 50000

$global.Range("A11").Value = "Green waste C/N"
$global.Range("B11").Value = 10

$global.Range("A12").Value = "Weight import"
$global.Range("B12").Value = 0

$global.Range("A14").Select()

# "global" becomes the active/selected sheet in the saved workbook.
$global.Activate()
